# Update the "assets" inspection sheet:
#  - R2 (row 3) no longer uses the shared "#" comment marker in column B
#  - R1 / R2 (rows 2 and 3) now specify an explicit port (column E) of 22
#  - the active selection on the sheet moves from E16 to H8

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("assets")

# Remove the "#" comment placeholder from B3
$ws.Range("B3").ClearContents()

# Record port 22 for the first two asset rows
$ws.Range("E2").Value2 = 22
$ws.Range("E3").Value2 = 22

# Move the saved selection to H8 (bottom-right frozen pane)
$ws.Activate() | Out-Null
$ws.Range("H8").Select() | Out-Null
